# Auto-generated edit script: update crypto price/volume table to latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.140.47"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "'1.851.85"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'237.34"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").Value = "'0.6883"
$ws.Range("E6").Value = "  -5.20%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.07747"
$ws.Range("E8").Value = "  +8.60%  "
$ws.Range("D9").Value = "'0.3037"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").Value = "'23.19"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").Value = "'0.08159"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'1.870.52"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "'0.7233"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "'5.197"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "'89.03"
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").Value = "'29.136.39"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000007826"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'5.731"
$ws.Range("E18").Value = "  -4.31%  "
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "'234.19"
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'2.100.64"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'7.490"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").Value = "'161.73"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "'8.957"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("E27").Value = "  -6.69%  "
$ws.Range("D28").Value = "'18.06"
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "'1.962"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "'1.400"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").Value = "'4.519"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "'1.482"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "'4.001"
$ws.Range("E33").Value = "  -4.27%  "
$ws.Range("D34").Value = "'0.05192"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").Value = "'1.178"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("D36").Value = "'0.7031"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("D37").Value = "'1.024"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").Value = "'2.651"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").Value = "'0.01849"
$ws.Range("D40").Value = "'2.675"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "'0.9094"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "'1.092.74"
$ws.Range("E42").Value = "  +5.30%  "
$ws.Range("D43").Value = "'5.991"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'0.4275"
$ws.Range("E44").Value = "  -4.30%  "
$ws.Range("D45").Value = "'70.46"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "'102.40"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").Value = "'1.754"
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("D49").Value = "'1.997.32"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").Value = "'9.138"
$ws.Range("E50").Value = "  -3.96%  "
$ws.Range("D51").Value = "'6.915"
$ws.Range("E51").Value = "  -7.35%  "
